$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression: only B2 changes
$ws.Range("B2").Value = 14975636417100490

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.01960181934381909
$ws.Range("C3").Value = 0.01988817470148027
$ws.Range("D3").Value = 2770504191989588

# Row 4 - was GradientBoostingRegressor, now DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.02185858952465563
$ws.Range("C4").Value = 0.021801488983171
$ws.Range("D4").Value = 126970743358660

# Row 5 - was AdaBoostRegressor, now MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 50566869823150.85
$ws.Range("C5").Value = 19624965754618.2
$ws.Range("D5").Value = 583084206126640
